$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4672152491505415
$ws.Range("C2").Value = 0.04465201330235402
$ws.Range("E2").Value = 0.537108540994808
$ws.Range("F2").Value = 2.560859243357569
$ws.Range("G2").Value = 0.84491862659074
$ws.Range("H2").Value = 0.9063570786317854
$ws.Range("J2").Value = 0.06869674885786026
$ws.Range("K2").Value = 0.4350942793186618
$ws.Range("N2").Value = 1.626256284582072
$ws.Range("B3").Value = 0.4274980773417099
$ws.Range("C3").Value = 0.03890564994279089
$ws.Range("E3").Value = 0.5144342463180038
$ws.Range("F3").Value = 2.526115842956088
$ws.Range("G3").Value = 0.8459506997002819
$ws.Range("H3").Value = 0.9117455477261416
$ws.Range("J3").Value = 0.06898881152527991
$ws.Range("K3").Value = 0.3922885727566552
$ws.Range("N3").Value = 1.646021725941028
$ws.Range("B4").Value = 0.4032828181369155
$ws.Range("C4").Value = 0.03536863920284361
$ws.Range("E4").Value = 0.500781382872745
$ws.Range("F4").Value = 2.506312281921367
$ws.Range("G4").Value = 0.847173483273238
$ws.Range("H4").Value = 0.9154916115074059
$ws.Range("J4").Value = 0.06921040242070831
$ws.Range("K4").Value = 0.3661267193897402
$ws.Range("N4").Value = 1.658768168530665
$ws.Range("B5").Value = 0.3934581660613787
$ws.Range("C5").Value = 0.03392508530556881
$ws.Range("E5").Value = 0.4952853012332383
$ws.Range("F5").Value = 2.498625926132618
$ws.Range("G5").Value = 0.8478195481917226
$ws.Range("H5").Value = 0.9171281187653193
$ws.Range("J5").Value = 0.06931130819503295
$ws.Range("K5").Value = 0.3554961473354297
$ws.Range("N5").Value = 1.664115819386534
$ws.Range("B6").Value = 0.3918294090601933
$ws.Range("C6").Value = 0.03368525181383575
$ws.Range("E6").Value = 0.49437675839836
$ws.Range("F6").Value = 2.497372769405558
$ws.Range("G6").Value = 0.8479357412677331
$ws.Range("H6").Value = 0.9174064988107915
$ws.Range("J6").Value = 0.06932870344663655
$ws.Range("K6").Value = 0.3537328017572179
$ws.Range("N6").Value = 1.665013052766056
$ws.Range("B7").Value = 0.4031501440055933
$ws.Range("C7").Value = 0.03534917981103547
$ws.Range("E7").Value = 0.5007069873744427
$ws.Range("F7").Value = 2.506207068197938
$ws.Range("G7").Value = 0.8471815984927531
$ws.Range("H7").Value = 0.9155132368957908
$ws.Range("J7").Value = 0.06921172035766077
$ws.Range("K7").Value = 0.3659832276783561
$ws.Range("N7").Value = 1.658839667977778
$ws.Range("B8").Value = 0.45348533383671
$ws.Range("C8").Value = 0.04267247501186944
$ws.Range("E8").Value = 0.5292344735509147
$ws.Range("F8").Value = 2.548561910007706
$ws.Range("G8").Value = 0.8451520012447133
$ws.Range("H8").Value = 0.9081241763912971
$ws.Range("J8").Value = 0.06878866709404008
$ws.Range("K8").Value = 0.4203098368181486
$ws.Range("N8").Value = 1.63294465753568
$ws.Range("B9").Value = 0.55354782724919
$ws.Range("C9").Value = 0.05696498608108413
$ws.Range("E9").Value = 0.5873226712328119
$ws.Range("F9").Value = 2.643792738976884
$ws.Range("G9").Value = 0.8458634775144986
$ws.Range("H9").Value = 0.8971090591974473
$ws.Range("J9").Value = 0.06829544232022045
$ws.Range("K9").Value = 0.5278029892918141
$ws.Range("N9").Value = 1.587012004902398
$ws.Range("B10").Value = 0.6278934740698503
$ws.Range("C10").Value = 0.06742603684119786
$ws.Range("E10").Value = 0.6313266061686278
$ws.Range("F10").Value = 2.721248199289533
$ws.Range("G10").Value = 0.8492719617082258
$ws.Range("H10").Value = 0.8911397569834065
$ws.Range("J10").Value = 0.0681396666176397
$ws.Range("K10").Value = 0.6073693197584191
$ws.Range("N10").Value = 1.556225101995263
$ws.Range("B11").Value = 0.6618966803484909
$ws.Range("C11").Value = 0.07217694238701711
$ws.Range("E11").Value = 0.6516374888983876
$ws.Range("F11").Value = 2.758126517966758
$ws.Range("G11").Value = 0.8514548540032933
$ws.Range("H11").Value = 0.8888864585131273
$ws.Range("J11").Value = 0.06811399923369166
$ws.Range("K11").Value = 0.6436964247957633
$ws.Range("N11").Value = 1.542863076085673
$ws.Range("B12").Value = 0.6747990514069215
$ws.Range("C12").Value = 0.07397487269162184
$ws.Range("E12").Value = 0.659371108588175
$ws.Range("F12").Value = 2.772328749703462
$ws.Range("G12").Value = 0.8523728255790246
$ws.Range("H12").Value = 0.888099746567022
$ws.Range("J12").Value = 0.06811080537648806
$ws.Range("K12").Value = 0.657471507123546
$ws.Range("N12").Value = 1.537895835471433
$ws.Range("B13").Value = 0.6720191378553579
$ws.Range("C13").Value = 0.07358770670845161
$ws.Range("E13").Value = 0.6577036500733868
$ws.Range("F13").Value = 2.769259480627397
$ws.Range("G13").Value = 0.8521710540191805
$ws.Range("H13").Value = 0.888266216961469
$ws.Range("J13").Value = 0.06811120256532277
$ws.Range("K13").Value = 0.6545039636006607
$ws.Range("N13").Value = 1.538961495032197
$ws.Range("B14").Value = 0.6629576444549343
$ws.Range("C14").Value = 0.07232488193980657
$ws.Range("E14").Value = 0.6522728891917922
$ws.Range("F14").Value = 2.759290184897196
$ws.Range("G14").Value = 0.8515285426888539
$ws.Range("H14").Value = 0.8888204006179024
$ws.Range("J14").Value = 0.06811360553215096
$ws.Range("K14").Value = 0.6448293334467792
$ws.Range("N14").Value = 1.542452559873869
$ws.Range("B15").Value = 0.6574106099192818
$ws.Range("C15").Value = 0.0715512175205788
$ws.Range("E15").Value = 0.6489519085770183
$ws.Range("F15").Value = 2.753214622051473
$ws.Range("G15").Value = 0.8511468960758179
$ws.Range("H15").Value = 0.8891685257351867
$ws.Range("J15").Value = 0.06811592807060407
$ws.Range("K15").Value = 0.6389057855454325
$ws.Range("N15").Value = 1.544603011878404
$ws.Range("B16").Value = 0.6256749481014481
$ws.Range("C16").Value = 0.06711539530499522
$ws.Range("E16").Value = 0.6300051569403706
$ws.Range("F16").Value = 2.718871260336385
$ws.Range("G16").Value = 0.8491420655161193
$ws.Range("H16").Value = 0.891296324189824
$ws.Range("J16").Value = 0.06814225577548072
$ws.Range("K16").Value = 0.6049978981441768
$ws.Range("N16").Value = 1.557111298836647
$ws.Range("B17").Value = 0.6062528604953172
$ws.Range("C17").Value = 0.06439214282241323
$ws.Range("E17").Value = 0.6184571630959823
$ws.Range("F17").Value = 2.698224273509084
$ws.Range("G17").Value = 0.8480744306005192
$ws.Range("H17").Value = 0.8927201053080722
$ws.Range("J17").Value = 0.0681700016390856
$ws.Range("K17").Value = 0.584230180434588
$ws.Range("N17").Value = 1.564949574621941
$ws.Range("B18").Value = 0.595099027124121
$ws.Range("C18").Value = 0.06282505848400888
$ws.Range("E18").Value = 0.6118426570160267
$ws.Range("F18").Value = 2.686503307477437
$ws.Range("G18").Value = 0.8475198499634189
$ws.Range("H18").Value = 0.8935825215451558
$ws.Range("J18").Value = 0.06819021275512327
$ws.Range("K18").Value = 0.5722975720593979
$ws.Range("N18").Value = 1.569518462442002
$ws.Range("B19").Value = 0.5913254981339264
$ws.Range("C19").Value = 0.06229434334782979
$ws.Range("E19").Value = 0.6096078335726958
$ws.Range("F19").Value = 2.682561321562559
$ws.Range("G19").Value = 0.8473422838343083
$ws.Range("H19").Value = 0.8938819880300883
$ws.Range("J19").Value = 0.06819778543949795
$ws.Range("K19").Value = 0.5682595369829926
$ws.Range("N19").Value = 1.571075797733329
$ws.Range("B20").Value = 0.6083185924396446
$ws.Range("C20").Value = 0.06468211445834982
$ws.Range("E20").Value = 0.6196836096984129
$ws.Range("F20").Value = 2.700406169985257
$ws.Range("G20").Value = 0.8481819215980693
$ws.Range("H20").Value = 0.892564039054534
$ws.Range("J20").Value = 0.0681666077413503
$ws.Range("K20").Value = 0.5864396539442112
$ws.Range("N20").Value = 1.564108911912938
$ws.Range("B21").Value = 0.6656185182588104
$ws.Range("C21").Value = 0.07269583498512588
$ws.Range("E21").Value = 0.6538668852341658
$ws.Range("F21").Value = 2.762211963694313
$ws.Range("G21").Value = 0.8517147810184724
$ws.Range("H21").Value = 0.8886558160965166
$ws.Range("J21").Value = 0.06811272239984945
$ws.Range("K21").Value = 0.647670498164814
$ws.Range("N21").Value = 1.541424632946592
$ws.Range("B22").Value = 0.703219332055653
$ws.Range("C22").Value = 0.07792665140526367
$ws.Range("E22").Value = 0.6764544343770638
$ws.Range("F22").Value = 2.803988700391557
$ws.Range("G22").Value = 0.8545564101248857
$ws.Range("H22").Value = 0.8864895978172882
$ws.Range("J22").Value = 0.06811555526447322
$ws.Range("K22").Value = 0.6877978692240276
$ws.Range("N22").Value = 1.527139381433297
$ws.Range("B23").Value = 0.6831372342690827
$ws.Range("C23").Value = 0.07513547172632684
$ws.Range("E23").Value = 0.6643763990048086
$ws.Range("F23").Value = 2.781564821897263
$ws.Range("G23").Value = 0.852990895818678
$ws.Range("H23").Value = 0.8876102095059224
$ws.Range("J23").Value = 0.06811055284152445
$ws.Range("K23").Value = 0.6663711789002207
$ws.Range("N23").Value = 1.534714192519088
$ws.Range("B24").Value = 0.6073846370609886
$ws.Range("C24").Value = 0.06455102279922187
$ws.Range("E24").Value = 0.6191290566111292
$ws.Range("F24").Value = 2.699419269700684
$ws.Range("G24").Value = 0.8481331404868797
$ws.Range("H24").Value = 0.8926344599697273
$ws.Range("J24").Value = 0.06816812885697487
$ws.Range("K24").Value = 0.5854407290631798
$ws.Range("N24").Value = 1.564488780658971
$ws.Range("B25").Value = 0.5263328217486105
$ws.Range("C25").Value = 0.05310559872765452
$ws.Range("E25").Value = 0.5713766586307685
$ws.Range("F25").Value = 2.616720037370186
$ws.Range("G25").Value = 0.8451660980378506
$ws.Range("H25").Value = 0.8997164105117292
$ws.Range("J25").Value = 0.06839270344480752
$ws.Range("K25").Value = 0.4986199229420549
$ws.Range("N25").Value = 1.598918587214236
